# Weekly price update: a new "Jengibre" price record for Vega Modelo de
# Temuco is inserted at row 20 (pushing the existing rows 20-94 down to
# 21-95), then the new row is populated with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 20..94 down to 21..95, leaving a blank row 20 behind.
$ws.Rows.Item(20).Insert()

# Populate the new row 20 with the latest observation.
$ws.Cells.Item(20, 1).Value = 10
$ws.Cells.Item(20, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(20, 3).Value = "La Araucanía"
$ws.Cells.Item(20, 4).Value = 44414
$ws.Cells.Item(20, 5).Value = 9
$ws.Cells.Item(20, 6).Value = 100114007
$ws.Cells.Item(20, 7).Value = "Jengibre"
$ws.Cells.Item(20, 8).Value = "Sin especificar"
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 80
$ws.Cells.Item(20, 11).Value = 20000
$ws.Cells.Item(20, 12).Value = 20000
$ws.Cells.Item(20, 13).Value = 20000
$ws.Cells.Item(20, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(20, 15).Value = "Perú"
$ws.Cells.Item(20, 16).Value = 1538
$ws.Cells.Item(20, 17).Value = 13
$ws.Cells.Item(20, 18).Value = "Hortaliza"
